# "Final Changes For TechTalk"
# Slide 9 (last slide): add the yellow bold-italic GitHub link to the
# (previously empty) title placeholder, and reposition the picture.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)

# --- 1) Title placeholder: add the bold/italic yellow hyperlink text ---
$titleShape = $s.Shapes.Item(1)
$tr = $titleShape.TextFrame.TextRange

# Seed the (still-empty) range's language before typing so the new run
# keeps "lang=en-US" once the text is assigned.
$tr.LanguageID = "en-US"
$tr.Text = "https://github.com/ImSraone/Bootstrap"

$full = $titleShape.TextFrame.TextRange
$full.LanguageID = "en-US"
$full.Font.Bold = $true
$full.Font.Italic = $true
$full.Font.Color.RGB = 65535   # RGB(255,255,0) -> yellow (0x00FFFF00, BGR-packed)

# --- 2) Picture: move it to its new position ---
$pic = $s.Shapes.Item(2)
$pic.Left = 203.8749842519685   # 2589212 EMU
$pic.Top = 108                  # 1371600 EMU
